$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 80 (id 79) to the dataset
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "F"
$ws.Cells.Item(80, 3).Value = "P"
$ws.Cells.Item(80, 4).Value = "socialmente"
$ws.Cells.Item(80, 5).Value = 1.6
$ws.Cells.Item(80, 6).Value = 19

# Update the frozen-pane view / selection to reflect scrolled position
$excel.ActiveWindow.ScrollRow = 65
$ws.Range("E81").Select()
